$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.064.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.85%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.296.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.64%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.70%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.38%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.12%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.93%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.88%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.651.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.79%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.298.39"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.66%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.809"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.53%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.958.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.88%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0924"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.67%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.34"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.49%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.52%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.98%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.46%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.29%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.59%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.74%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.69"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.33%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.69%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.58%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.53%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.84%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.969.68"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.04"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.53%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.95"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +17.99%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.520.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.63%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.17%  "
